$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.742.39'
$ws.Range("E2").Value = '  -4.05%  '

$ws.Range("D3").Value = '1.817.65'
$ws.Range("E3").Value = '  -3.03%  '

$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").Value = "'278.27"
$ws.Range("E5").Value = '  -7.67%  '

$ws.Range("E6").Value = '  -0.15%  '

$ws.Range("D7").Value = "'0.5094"
$ws.Range("E7").Value = '  -5.06%  '

$ws.Range("D8").Value = "'0.3535"
$ws.Range("E8").Value = '  -5.63%  '

$ws.Range("D9").Value = "'44.60"
$ws.Range("E9").Value = '  -1.98%  '

$ws.Range("D10").Value = "'0.06665"
$ws.Range("E10").Value = '  -7.38%  '

$ws.Range("E11").Value = '  -6.96%  '

$ws.Range("D12").Value = "'0.8276"
$ws.Range("E12").Value = '  -7.01%  '

$ws.Range("E13").Value = '  -3.30%  '

$ws.Range("D14").Value = '1.792.78'
$ws.Range("E14").Value = '  -4.34%  '

$ws.Range("D15").Value = "'5.076"
$ws.Range("E15").Value = '  -4.56%  '

$ws.Range("D16").Value = "'87.79"
$ws.Range("E16").Value = '  -5.98%  '

$ws.Range("D17").Value = "'0.9999"
$ws.Range("E17").Value = '  -0.26%  '

$ws.Range("D18").Value = "'14.10"
$ws.Range("E18").Value = '  -5.01%  '

$ws.Range("D19").Value = "'0.000008032"
$ws.Range("E19").Value = '  -5.82%  '

$ws.Range("D21").Value = '25.782.52'
$ws.Range("E21").Value = '  -4.03%  '

$ws.Range("D22").Value = "'4.746"
$ws.Range("E22").Value = '  -4.91%  '

$ws.Range("D23").Value = "'10.00"
$ws.Range("E23").Value = '  -5.86%  '

$ws.Range("D24").Value = "'6.110"
$ws.Range("E24").Value = '  -4.58%  '

$ws.Range("D25").Value = "'2.227"
$ws.Range("E25").Value = '  -2.82%  '

$ws.Range("D26").Value = "'141.98"
$ws.Range("E26").Value = '  -2.96%  '

$ws.Range("E27").Value = '  -3.51%  '

$ws.Range("D28").Value = "'17.11"
$ws.Range("E28").Value = '  -5.52%  '

$ws.Range("D29").Value = "'109.36"
$ws.Range("E29").Value = '  -4.13%  '

$ws.Range("D30").Value = "'4.320"
$ws.Range("E30").Value = '  -8.44%  '

$ws.Range("E31").Value = '  -8.29%  '

$ws.Range("D32").Value = "'0.08766"
$ws.Range("E32").Value = '  -4.07%  '

$ws.Range("D33").Value = "'0.04891"
$ws.Range("E33").Value = '  -2.47%  '

$ws.Range("D34").Value = "'0.7296"
$ws.Range("E34").Value = '  -9.99%  '

$ws.Range("D35").Value = "'1.140"
$ws.Range("E35").Value = '  -3.02%  '

$ws.Range("D36").Value = "'2.866"
$ws.Range("E36").Value = '  -2.78%  '

$ws.Range("D37").Value = "'0.9991"
$ws.Range("E37").Value = '  -0.36%  '

$ws.Range("D38").Value = "'3.140"
$ws.Range("E38").Value = '  -2.68%  '

$ws.Range("D39").Value = "'2.381"
$ws.Range("E39").Value = '  -9.62%  '

$ws.Range("E40").Value = '  -5.29%  '

$ws.Range("E41").Value = '  -14.56%  '

$ws.Range("D42").Value = "'0.9657"
$ws.Range("E42").Value = '  -9.77%  '

$ws.Range("E43").Value = '  -6.12%  '

$ws.Range("D44").Value = "'111.29"
$ws.Range("E44").Value = '  -3.35%  '

$ws.Range("D45").Value = "'8.033"
$ws.Range("E45").Value = '  -9.61%  '

$ws.Range("E46").Value = '  -0.13%  '

$ws.Range("D47").Value = "'0.4559"
$ws.Range("E47").Value = '  -11.30%  '

$ws.Range("D48").Value = "'0.1367"
$ws.Range("E48").Value = '  -8.80%  '

$ws.Range("D49").Value = "'36.52"
$ws.Range("E49").Value = '  -2.96%  '

$ws.Range("D50").Value = "'9.199"
$ws.Range("E50").Value = '  -7.95%  '

$ws.Range("D51").Value = "'1.502"
$ws.Range("E51").Value = '  -8.79%  '
